$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the credentials table -----------------------------------
# Header row: EMAIL / PASSWORD (was USERNAME / PASSWORD)
$ws.Range("B1").Value = "PASSWORD"
$ws.Range("A1").Value = "EMAIL"

# Data rows: three mock email/password pairs replacing the old admin/password row
$ws.Range("A2").Value = "localhost@microsoft.com"
$ws.Range("A3").Value = "root@ubuntu.org"
$ws.Range("A4").Value = "admin@stackworks.online"

$ws.Range("B2").Value = "syspass"
$ws.Range("B3").Value = "sudologin"
$ws.Range("B4").Value = "letmein"

# --- Formatting: left/center align the whole table --------------------
$table = $ws.Range("A1:B4")
$table.HorizontalAlignment = -4131   # xlLeft
$table.VerticalAlignment = -4108     # xlCenter

# --- Column widths (approx. best-fit for the new, wider content) ------
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 10.333333333333334

# --- Selection / active cell ------------------------------------------
[void]$ws.Range("D12").Select()

# --- Page setup: paper size + orientation ------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9    # xlPaperA4
$ps.Orientation = 1  # xlPortrait
